$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H135").Value = 2229.4546
$ws.Range("I135").Value = 2352.4
$ws.Range("J135").Value = 1000
$ws.Range("K135").Value = 21171.6
$ws.Range("L135").Value = 9000
$ws.Range("M135").Value = -18636.6
$ws.Range("N135").Value = -14070

$ws.Range("H137").Value = 1050
$ws.Range("I137").Value = 862.125
$ws.Range("K137").Value = 2586.375
$ws.Range("M137").Value = -36.375

$ws.Range("H138").Value = 2390.28
$ws.Range("I138").Value = 875.8043
$ws.Range("J138").Value = 3680.389
$ws.Range("K138").Value = 2627.4129
$ws.Range("L138").Value = 11041.167
$ws.Range("M138").Value = 2512.5871
$ws.Range("N138").Value = -21321.167

$ws.Range("H139").Value = 113320
$ws.Range("J139").Value = 113320
$ws.Range("L139").Value = 113320
$ws.Range("N139").Value = -123600

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 20624.809
$ws.Range("I2").Value = 23536.25
$ws.Range("J2").Value = 4611.875
$ws.Range("K2").Value = 23536.25
$ws.Range("L2").Value = 4611.875
$ws.Range("M2").Value = -23423.25
$ws.Range("N2").Value = -4837.875

$ws.Range("H32").Value = 1454.46
$ws.Range("I32").Value = 1454.46
$ws.Range("K32").Value = 1454.46
$ws.Range("M32").Value = -1167.46

$ws.Range("H74").Value = 1009.86487
$ws.Range("I74").Value = 1009.86487
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 1009.86487
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -135.86487
$ws.Range("N74").ClearContents()

$ws.Range("H77").Value = 1009.86487
$ws.Range("I77").Value = 1009.86487
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 5049.32435
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -681.3243499999999
$ws.Range("N77").ClearContents()

$ws.Range("H111").Value = 30000
$ws.Range("J111").Value = 30000
$ws.Range("L111").Value = 30000
$ws.Range("N111").Value = -38180

$ws.Range("H116").Value = 20624.809
$ws.Range("I116").Value = 23536.25
$ws.Range("J116").Value = 4611.875
$ws.Range("K116").Value = 23536.25
$ws.Range("L116").Value = 4611.875
$ws.Range("M116").Value = -21242.25
$ws.Range("N116").Value = -9199.875

$ws.Range("H121").Value = 13556.714
$ws.Range("J121").Value = 13556.714
$ws.Range("L121").Value = 13556.714
$ws.Range("N121").Value = -17050.714

$ws.Range("H132").Value = 1562.7894
$ws.Range("I132").Value = 1253.7576
$ws.Range("K132").Value = 3761.2728
$ws.Range("M132").Value = -1231.2728

$ws.Range("H139").Value = 49715
$ws.Range("J139").Value = 49715
$ws.Range("L139").Value = 49715
$ws.Range("N139").Value = -59995

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 20624.809
$ws.Range("I3").Value = 23536.25
$ws.Range("J3").Value = 4611.875
$ws.Range("K3").Value = 23536.25
$ws.Range("L3").Value = 4611.875
$ws.Range("M3").Value = -23422.25
$ws.Range("N3").Value = -4839.875

$ws.Range("H86").Value = 2458.4443
$ws.Range("I86").Value = 1944.9166
$ws.Range("J86").Value = 6566.6665
$ws.Range("K86").Value = 1944.9166
$ws.Range("L86").Value = 6566.6665
$ws.Range("M86").Value = -821.9166
$ws.Range("N86").Value = -8812.666499999999

$ws.Range("H89").Value = 2458.4443
$ws.Range("I89").Value = 1944.9166
$ws.Range("J89").Value = 6566.6665
$ws.Range("K89").Value = 9724.583000000001
$ws.Range("L89").Value = 32833.3325
$ws.Range("M89").Value = -4108.583000000001
$ws.Range("N89").Value = -44065.3325

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 4255.3516
$ws.Range("I58").Value = 1568.6666
$ws.Range("J58").Value = 9215.385
$ws.Range("K58").Value = 1568.6666
$ws.Range("L58").Value = 9215.385
$ws.Range("M58").Value = -1365.6666
$ws.Range("N58").Value = -9621.385

$ws.Range("H60").Value = 8247
$ws.Range("J60").Value = 8247
$ws.Range("L60").Value = 8247
$ws.Range("N60").Value = -9269

$ws.Range("H64").Value = 29800
$ws.Range("J64").Value = 29800
$ws.Range("L64").Value = 29800
$ws.Range("N64").Value = -30296

$ws.Range("H67").Value = 29800
$ws.Range("J67").Value = 29800
$ws.Range("L67").Value = 29800
$ws.Range("N67").Value = -31516

$ws.Range("H74").Value = 11698.375
$ws.Range("J74").Value = 13681.167
$ws.Range("L74").Value = 13681.167
$ws.Range("N74").Value = -15429.167

$ws.Range("H77").Value = 11698.375
$ws.Range("J77").Value = 13681.167
$ws.Range("L77").Value = 41043.501
$ws.Range("N77").Value = -49779.501

$ws.Range("H98").Value = 35000
$ws.Range("J98").Value = 35000
$ws.Range("L98").Value = 35000
$ws.Range("N98").Value = -39492

$ws.Range("H136").Value = 4255.3516
$ws.Range("I136").Value = 1568.6666
$ws.Range("J136").Value = 9215.385
$ws.Range("K136").Value = 4705.9998
$ws.Range("L136").Value = 27646.155
$ws.Range("M136").Value = -2155.9998
$ws.Range("N136").Value = -32746.155

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H121").Value = 959.76
$ws.Range("J121").Value = 970.5833
$ws.Range("L121").Value = 2911.7499
$ws.Range("N121").Value = -5531.7499

$ws.Range("H123").Value = 8449.5
$ws.Range("I123").Value = 1200
$ws.Range("J123").Value = 9899.4
$ws.Range("K123").Value = 3600
$ws.Range("L123").Value = 29698.2
$ws.Range("M123").Value = -1150
$ws.Range("N123").Value = -34598.2

$ws.Range("H134").Value = 4347.7407
$ws.Range("I134").Value = 2522.8823
$ws.Range("J134").Value = 7450
$ws.Range("K134").Value = 7568.646900000001
$ws.Range("L134").Value = 22350
$ws.Range("M134").Value = -2498.646900000001
$ws.Range("N134").Value = -32490

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H42").Value = 24500
$ws.Range("J42").Value = 24500
$ws.Range("L42").Value = 24500
$ws.Range("N42").Value = -25470

$ws.Range("H110").Value = 29916.666
$ws.Range("J110").Value = 29916.666
$ws.Range("L110").Value = 29916.666
$ws.Range("N110").Value = -38096.666

$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

$ws.Range("H113").Value = 4227.5557
$ws.Range("I113").Value = 4638
$ws.Range("J113").Value = 3406.6667
$ws.Range("K113").Value = 4638
$ws.Range("L113").Value = 3406.6667
$ws.Range("M113").Value = -2468
$ws.Range("N113").Value = -7746.6667

$ws.Range("H115").Value = 24500
$ws.Range("J115").Value = 24500
$ws.Range("L115").Value = 24500
$ws.Range("N115").Value = -26850

$ws.Range("H116").Value = 33500
$ws.Range("J116").Value = 33500
$ws.Range("L116").Value = 33500
$ws.Range("N116").Value = -42678

$ws.Range("H118").Value = 20000
$ws.Range("J118").Value = 20000
$ws.Range("L118").Value = 20000
$ws.Range("N118").Value = -23314

$ws.Range("H119").Value = 26593.9
$ws.Range("J119").Value = 26593.9
$ws.Range("L119").Value = 26593.9
$ws.Range("N119").Value = -36269.9

$ws.Range("H132").Value = 2766.353
$ws.Range("I132").Value = 2420.7036
$ws.Range("K132").Value = 7262.110799999999
$ws.Range("M132").Value = -4732.110799999999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 3464.6667
$ws.Range("I61").Value = 3464.6667
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 3464.6667
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -3262.6667
$ws.Range("N61").ClearContents()

$ws.Range("H113").Value = 3464.6667
$ws.Range("I113").Value = 3464.6667
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 3464.6667
$ws.Range("L113").Value = 3464.6667
$ws.Range("M113").Value = -1294.6667
$ws.Range("N113").ClearContents()
